$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$bullet = [char]0x2022

# Row 5 - month heading, was "Shahrivar, Mehr o Aban 99", now "Azar 99"
$ws.Cells.Item(5, 1).Value = "آذر 99"

# Row 6 - was "* Code Refactoring" / 6 / "Adjusted orientation widget for the phantom"
$ws.Cells.Item(6, 2).Value = "* Documentation"
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 5).Value = "$bullet Fixed tracker's crash"

# Row 7 - "* GUI" stays, hours 8 -> 5, task text changes
$ws.Cells.Item(7, 2).Value = "* GUI"
$ws.Cells.Item(7, 3).Value = 5
$ws.Cells.Item(7, 5).Value = "$bullet Added tracker status splash messages"

# Row 8 - "* Registration" stays, hours 14 -> 6, task text changes
$ws.Cells.Item(8, 2).Value = "* Registration"
$ws.Cells.Item(8, 3).Value = 6
$ws.Cells.Item(8, 5).Value = "$bullet Complete and working registration process"

# Row 9 - was "* Tracker" / 2 / "Record/Load Tracker Centerline"
$ws.Cells.Item(9, 2).Value = "* Tracker"
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(9, 5).Value = "$bullet Presentations @IACT"

# Row 10 - was "* 2D/3D Views" / 2 / "Integrated Registration Process"
$ws.Cells.Item(10, 2).Value = "* 2D/3D Views"
$ws.Cells.Item(10, 3).Value = 4
$ws.Cells.Item(10, 5).Value = "$bullet Document of the code"

# Row 11 - was "* Patients Database" / 1
$ws.Cells.Item(11, 2).Value = "* Patients / Database"
$ws.Cells.Item(11, 3).Value = 1

# Row 12 - was "* Meetings & other" / 2
$ws.Cells.Item(12, 2).Value = "* Meetings & Presentations"
$ws.Cells.Item(12, 3).Value = 12

# Row 14 - @Home hours 0 -> 6 (D15 = C13-D14 will recalc automatically)
$ws.Cells.Item(14, 4).Value = 6

$wb.Save()
